# Update cryptos list (Price / Volume(1h) columns) with refreshed figures.
# Leading "'" on some D-column values forces text storage (matches the
# original inlineStr "Text" cell type) so numeric-looking strings like
# "0.6279" aren't reinterpreted as Excel numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.412.93"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.849.83"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D5").Value = "'240.50"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.6279"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.07637"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'24.74"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'5.035"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'0.6792"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'0.00001069"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "'83.28"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'6.167"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "29.438.85"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'226.79"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'12.33"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D21").Value = "'7.469"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'157.79"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "'0.1380"
$ws.Range("D25").Value = "'8.419"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'17.69"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'1.389"
$ws.Range("E27").Value = "  +7.18%  "
$ws.Range("D28").Value = "'1.465"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "'0.05588"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "'4.128"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "'4.060"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "'1.838"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "'1.164"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'0.6957"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "'2.589"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "1.230.09"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'2.725"
$ws.Range("D39").Value = "'6.416"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'0.9045"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'66.04"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'7.172"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").Value = "'0.4013"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "'8.984"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'0.1145"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").Value = "'0.05704"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'0.4634"
$ws.Range("E51").Value = "  +0.21%  "
